$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 50,4
$arr[0,0] = "Bitcoin"
$arr[0,1] = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$arr[0,2] = "30.344.26"
$arr[0,3] = "  +0.20%  "
$arr[1,0] = "Ethereum"
$arr[1,1] = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$arr[1,2] = "1.873.74"
$arr[1,3] = "  +0.33%  "
$arr[2,0] = "TetherUSD"
$arr[2,1] = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$arr[2,2] = "0.9999"
$arr[2,3] = "  -0.07%  "
$arr[3,0] = "BNB"
$arr[3,1] = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$arr[3,2] = "245.25"
$arr[3,3] = "  +4.50%  "
$arr[4,0] = "USDC"
$arr[4,1] = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$arr[4,2] = "0.9998"
$arr[4,3] = "  -0.06%  "
$arr[5,0] = "XRP"
$arr[5,1] = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$arr[5,2] = "0.4738"
$arr[5,3] = "  +0.82%  "
$arr[6,0] = "Cardano"
$arr[6,1] = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$arr[6,2] = "0.2887"
$arr[6,3] = "  +1.11%  "
$arr[7,0] = "OKB"
$arr[7,1] = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$arr[7,2] = "42.86"
$arr[7,3] = "  +2.92%  "
$arr[8,0] = "Dogecoin"
$arr[8,1] = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$arr[8,2] = "0.06476"
$arr[8,3] = "  -1.36%  "
$arr[9,0] = "Solana"
$arr[9,1] = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$arr[9,2] = "21.15"
$arr[9,3] = "  -0.96%  "
$arr[10,0] = "TRON"
$arr[10,1] = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$arr[10,2] = "0.07779"
$arr[10,3] = "  -0.57%  "
$arr[11,0] = "WrappedEther"
$arr[11,1] = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$arr[11,2] = "1.871.42"
$arr[11,3] = "  +1.49%  "
$arr[12,0] = "Polygon"
$arr[12,1] = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$arr[12,2] = "0.7302"
$arr[12,3] = "  +5.16%  "
$arr[13,0] = "Litecoin"
$arr[13,1] = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$arr[13,2] = "95.47"
$arr[13,3] = "  -1.36%  "
$arr[14,0] = "Polkadot"
$arr[14,1] = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$arr[14,2] = "5.130"
$arr[14,3] = "  +0.94%  "
$arr[15,0] = "BitcoinCash"
$arr[15,1] = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$arr[15,2] = "275.39"
$arr[15,3] = "  +2.56%  "
$arr[16,0] = "WrappedBTC"
$arr[16,1] = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$arr[16,2] = "30.327.67"
$arr[16,3] = "  +0.86%  "
$arr[17,0] = "Avalanche"
$arr[17,1] = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$arr[17,2] = "13.38"
$arr[17,3] = "  -2.81%  "
$arr[18,0] = "ShibaInu"
$arr[18,1] = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$arr[18,2] = "0.000007541"
$arr[18,3] = "  -1.93%  "
$arr[19,0] = "Dai"
$arr[19,1] = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$arr[19,2] = "1.000"
$arr[19,3] = "  +0.01%  "
$arr[20,0] = "WrappedliquidstakedEther2.0"
$arr[20,1] = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$arr[20,2] = "2.125.49"
$arr[20,3] = "  +0.86%  "
$arr[21,0] = "BinanceUSD"
$arr[21,1] = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$arr[21,2] = "0.9998"
$arr[21,3] = "  -0.08%  "
$arr[22,0] = "Uniswap"
$arr[22,1] = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$arr[22,2] = "5.241"
$arr[22,3] = "  -0.16%  "
$arr[23,0] = "Chainlink"
$arr[23,1] = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$arr[23,2] = "6.162"
$arr[23,3] = "  +0.26%  "
$arr[24,0] = "Cosmos"
$arr[24,1] = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$arr[24,2] = "9.266"
$arr[24,3] = "  -3.16%  "
$arr[25,0] = "Monero"
$arr[25,1] = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$arr[25,2] = "164.73"
$arr[25,3] = "  -1.04%  "
$arr[26,0] = "EthereumClassic"
$arr[26,1] = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$arr[26,2] = "18.95"
$arr[26,3] = "  +0.47%  "
$arr[27,0] = "LidoDAOToken"
$arr[27,1] = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$arr[27,2] = "1.920"
$arr[27,3] = "  -0.87%  "
$arr[28,0] = "Toncoin"
$arr[28,1] = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$arr[28,2] = "1.380"
$arr[28,3] = "  +1.27%  "
$arr[29,0] = "Stellar"
$arr[29,1] = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$arr[29,2] = "0.09893"
$arr[29,3] = "  -0.07%  "
$arr[30,0] = "PancakeSwap"
$arr[30,1] = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$arr[30,2] = "1.518"
$arr[30,3] = "  +4.12%  "
$arr[31,0] = "Filecoin"
$arr[31,1] = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$arr[31,2] = "4.319"
$arr[31,3] = "  -0.80%  "
$arr[32,0] = "InternetComputer(DFINITY)"
$arr[32,1] = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$arr[32,2] = "4.052"
$arr[32,3] = "  +0.09%  "
$arr[33,0] = "Hedera"
$arr[33,1] = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$arr[33,2] = "0.04776"
$arr[33,3] = "  +1.03%  "
$arr[34,0] = "ARBITRUM"
$arr[34,1] = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$arr[34,2] = "1.123"
$arr[34,3] = "  -0.56%  "
$arr[35,0] = "ImmutableX"
$arr[35,1] = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$arr[35,2] = "0.6975"
$arr[35,3] = "  -0.73%  "
$arr[36,0] = "HuobiToken"
$arr[36,1] = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$arr[36,2] = "2.717"
$arr[36,3] = "  -0.02%  "
$arr[37,0] = "VeChain"
$arr[37,1] = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$arr[37,2] = "0.01846"
$arr[37,3] = "  -1.57%  "
$arr[38,0] = "MXToken"
$arr[38,1] = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$arr[38,2] = "2.750"
$arr[38,3] = "  -0.72%  "
$arr[39,0] = "FraxShare"
$arr[39,1] = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$arr[39,2] = "6.428"
$arr[39,3] = "  +1.69%  "
$arr[40,0] = "Aave"
$arr[40,1] = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$arr[40,2] = "69.94"
$arr[40,3] = "  -4.06%  "
$arr[41,0] = "TrustWalletToken"
$arr[41,1] = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$arr[41,2] = "0.8443"
$arr[41,3] = "  +1.05%  "
$arr[42,0] = "RenderToken"
$arr[42,1] = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$arr[42,2] = "1.914"
$arr[42,3] = "  -1.80%  "
$arr[43,0] = "PaxDollar"
$arr[43,1] = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$arr[43,2] = "0.9996"
$arr[43,3] = "  -0.09%  "
$arr[44,0] = "TheSandbox"
$arr[44,1] = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$arr[44,2] = "0.4136"
$arr[44,3] = "  -0.73%  "
$arr[45,0] = "Quant"
$arr[45,1] = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$arr[45,2] = "102.53"
$arr[45,3] = "  -0.55%  "
$arr[46,0] = "EnergySwap"
$arr[46,1] = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$arr[46,2] = "9.355"
$arr[46,3] = "  +2.13%  "
$arr[47,0] = "Aptos"
$arr[47,1] = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$arr[47,2] = "7.089"
$arr[47,3] = "  -0.27%  "
$arr[48,0] = "Maker"
$arr[48,1] = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$arr[48,2] = "927.19"
$arr[48,3] = "  -5.00%  "
$arr[49,0] = "Elrond"
$arr[49,1] = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$arr[49,2] = "35.29"
$arr[49,3] = "  +2.19%  "
$ws.Range("B2:E51").Value = $arr
